# Update to include OLS on spreadsheet
# Appends three new rows (1 / 10 / 100 Lin Regs - ordinary least squares
# models) below the existing RMI/B-Tree comparison table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Model names (column A), entered 10/1/100 so the row for "10 Lin Regs"
# landed first, then "1 Lin Reg" was inserted above it, then "100 Lin
# Regs" was added last.
$ws.Range("A8").Value = "10 Lin Regs"
$ws.Range("A7").Value = "1 Lin Reg"
$ws.Range("A9").Value = "100 Lin Regs"

# Build time column (B), filled top to bottom.
$ws.Range("B7").Value = "0.00729 s"
$ws.Range("B8").Value = "0.00838 s"
$ws.Range("B9").Value = "0.0297 s"

# Query time column (D), filled top to bottom.
$ws.Range("D7").Value = "0.0365 ms"
$ws.Range("D8").Value = "0.0433 ms"
$ws.Range("D9").Value = "0.0763 ms"

# Array accesses column (F) - plain numbers, not strings.
$ws.Range("F7").Value = 14.4
$ws.Range("F8").Value = 14.43
$ws.Range("F9").Value = 6.39

# Storage column (H), filled top to bottom.
$ws.Range("H7").Value = "0.0156 kB"
$ws.Range("H8").Value = "0.234 kB"
$ws.Range("H9").Value = "2.344 kB"

# Center the value columns, matching the rest of the table (each area has
# to be set individually - multi-area ranges only format the first area).
$ws.Range("B7:B9").HorizontalAlignment = -4108
$ws.Range("D7:D9").HorizontalAlignment = -4108
$ws.Range("F7:F9").HorizontalAlignment = -4108
$ws.Range("H7:H9").HorizontalAlignment = -4108

# Row 6 ("B-Tree") used to be the last row of the table and had no bottom
# border; now that more rows follow it, give it the same separator line
# the rest of the table uses.
$ws.Range("A6:I6").Borders.Item(9).LineStyle = 1

# Close off the bottom of the table under the new last row.
$ws.Range("A9:I9").Borders.Item(9).LineStyle = 1

$ws.Range("G13").Select()
